# Apply commit "add user list to project":
#  1. Fix ordering of the cfop list for Romit Roy Choudhury on the "PI hours" sheet.
#  2. Add a new "users" column to the "project hours" sheet, listing the PI(s)
#     associated with each project.

$wb = $excel.ActiveWorkbook

# --- 1. PI hours sheet: reorder cfop list for the Choudhury/RRC row -------------
$wsPI = $wb.Worksheets.Item("PI hours")
$wsPI.Range("G8").Value = "['cfop_RRC', 'cfop_CHOUDHURY']"

# --- 2. project hours sheet: add "users" column ---------------------------------
$wsProj = $wb.Worksheets.Item("project hours")

# Header cell, matching the style of the existing header row (B1:D1)
$wsProj.Range("E1").Value = "users"
$wsProj.Range("D1").Copy()
$wsProj.Range("E1").PasteSpecial(-4122) # xlPasteFormats

$users = @(
    "['Jonathan Hoff']",
    "['Alexander Hilll', 'Alexander Hill']",
    "['Hamidreza Jafarnejadsani']",
    "['Harshal Maske']",
    "['Nicole Chan']",
    "['Karun Koppula']",
    "['Mahanth Gowda']"
)

for ($i = 0; $i -lt $users.Length; $i++) {
    $row = $i + 2
    $wsProj.Cells.Item($row, 5).Value = $users[$i]
}
